$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 33.36960033333333
$ws.Range("H2").Value = 100.108801
$ws.Range("I2").Value = 0.07727383968381614
$ws.Range("J2").Value = 0.07727383968381614
$ws.Range("M2").Value = 38.45264233333334
$ws.Range("N2").Value = 115.357927
$ws.Range("O2").Value = 0.2975040117664333
$ws.Range("P2").Value = 0.2975040117664332
$ws.Range("Q2").Value = 1283.149306423948
$ws.Range("R2").Value = 11548.34375781553
$ws.Range("S2").Value = 0.02298927731053152
$ws.Range("T2").Value = 0.02298927731053151
$ws.Range("G3").Value = 33.36960033333333
$ws.Range("H3").Value = 100.108801
$ws.Range("I3").Value = 0.07727383968381614
$ws.Range("J3").Value = 0.07727383968381614
$ws.Range("O3").Value = 0.3694391181876273
$ws.Range("P3").Value = 0.3694391181876272
$ws.Range("Q3").Value = 1593.408927340771
$ws.Range("R3").Value = 14340.68034606695
$ws.Range("S3").Value = 0.02854797919176111
$ws.Range("T3").Value = 0.02854797919176111
$ws.Range("G4").Value = 33.36960033333333
$ws.Range("H4").Value = 100.108801
$ws.Range("I4").Value = 0.07727383968381614
$ws.Range("J4").Value = 0.07727383968381614
$ws.Range("M4").Value = 18.63107466666667
$ws.Range("N4").Value = 55.893224
$ws.Range("O4").Value = 0.1441466469015163
$ws.Range("P4").Value = 0.1441466469015162
$ws.Range("Q4").Value = 621.7115154071582
$ws.Range("R4").Value = 5595.403638664425
$ws.Range("S4").Value = 0.01113876488362742
$ws.Range("T4").Value = 0.01113876488362742
$ws.Range("G5").Value = 33.36960033333333
$ws.Range("H5").Value = 100.108801
$ws.Range("I5").Value = 0.07727383968381614
$ws.Range("J5").Value = 0.07727383968381614
$ws.Range("M5").Value = 24.41680433333333
$ws.Range("N5").Value = 73.25041299999999
$ws.Range("O5").Value = 0.1889102231444233
$ws.Range("P5").Value = 0.1889102231444233
$ws.Range("Q5").Value = 814.7790020205347
$ws.Range("R5").Value = 7333.011018184812
$ws.Range("S5").Value = 0.0145978182978961
$ws.Range("T5").Value = 0.0145978182978961
$ws.Range("I6").Value = 0.2551852590901843
$ws.Range("J6").Value = 0.2551852590901843
$ws.Range("M6").Value = 38.45264233333334
$ws.Range("N6").Value = 115.357927
$ws.Range("O6").Value = 0.2975040117664333
$ws.Range("P6").Value = 0.2975040117664332
$ws.Range("Q6").Value = 4237.408022572521
$ws.Range("R6").Value = 38136.67220315269
$ws.Range("S6").Value = 0.07591863832298651
$ws.Range("T6").Value = 0.0759186383229865
$ws.Range("I7").Value = 0.2551852590901843
$ws.Range("J7").Value = 0.2551852590901843
$ws.Range("O7").Value = 0.3694391181876273
$ws.Range("P7").Value = 0.3694391181876272
$ws.Range("S7").Value = 0.09427541709275888
$ws.Range("T7").Value = 0.09427541709275887
$ws.Range("I8").Value = 0.2551852590901843
$ws.Range("J8").Value = 0.2551852590901843
$ws.Range("M8").Value = 18.63107466666667
$ws.Range("N8").Value = 55.893224
$ws.Range("O8").Value = 0.1441466469015163
$ws.Range("P8").Value = 0.1441466469015162
$ws.Range("Q8").Value = 2053.108979541934
$ws.Range("R8").Value = 18477.98081587741
$ws.Range("S8").Value = 0.03678409943654474
$ws.Range("T8").Value = 0.03678409943654473
$ws.Range("I9").Value = 0.2551852590901843
$ws.Range("J9").Value = 0.2551852590901843
$ws.Range("M9").Value = 24.41680433333333
$ws.Range("N9").Value = 73.25041299999999
$ws.Range("O9").Value = 0.1889102231444233
$ws.Range("P9").Value = 0.1889102231444233
$ws.Range("Q9").Value = 2690.685380493622
$ws.Range("R9").Value = 24216.16842444259
$ws.Range("S9").Value = 0.04820710423789419
$ws.Range("T9").Value = 0.04820710423789418
$ws.Range("G10").Value = 13.90116633333333
$ws.Range("H10").Value = 41.703499
$ws.Range("I10").Value = 0.0321908709702775
$ws.Range("J10").Value = 0.0321908709702775
$ws.Range("M10").Value = 38.45264233333334
$ws.Range("N10").Value = 115.357927
$ws.Range("O10").Value = 0.2975040117664333
$ws.Range("P10").Value = 0.2975040117664332
$ws.Range("Q10").Value = 534.5365770318415
$ws.Range("R10").Value = 4810.829193286573
$ws.Range("S10").Value = 0.009576913255913172
$ws.Range("T10").Value = 0.00957691325591317
$ws.Range("G11").Value = 13.90116633333333
$ws.Range("H11").Value = 41.703499
$ws.Range("I11").Value = 0.0321908709702775
$ws.Range("J11").Value = 0.0321908709702775
$ws.Range("O11").Value = 0.3694391181876273
$ws.Range("P11").Value = 0.3694391181876272
$ws.Range("Q11").Value = 663.7850712840617
$ws.Range("R11").Value = 5974.065641556555
$ws.Range("S11").Value = 0.01189256698495101
$ws.Range("T11").Value = 0.01189256698495101
$ws.Range("G12").Value = 13.90116633333333
$ws.Range("H12").Value = 41.703499
$ws.Range("I12").Value = 0.0321908709702775
$ws.Range("J12").Value = 0.0321908709702775
$ws.Range("M12").Value = 18.63107466666667
$ws.Range("N12").Value = 55.893224
$ws.Range("O12").Value = 0.1441466469015163
$ws.Range("P12").Value = 0.1441466469015162
$ws.Range("Q12").Value = 258.9936679100862
$ws.Range("R12").Value = 2330.943011190776
$ws.Range("S12").Value = 0.004640206111204861
$ws.Range("T12").Value = 0.00464020611120486
$ws.Range("G13").Value = 13.90116633333333
$ws.Range("H13").Value = 41.703499
$ws.Range("I13").Value = 0.0321908709702775
$ws.Range("J13").Value = 0.0321908709702775
$ws.Range("M13").Value = 24.41680433333333
$ws.Range("N13").Value = 73.25041299999999
$ws.Range("O13").Value = 0.1889102231444233
$ws.Range("P13").Value = 0.1889102231444233
$ws.Range("Q13").Value = 339.4220583661208
$ws.Range("R13").Value = 3054.798525295087
$ws.Range("S13").Value = 0.006081184618208461
$ws.Range("T13").Value = 0.00608118461820846
$ws.Range("G14").Value = 274.366806
$ws.Range("H14").Value = 823.100418
$ws.Range("I14").Value = 0.635350030255722
$ws.Range("J14").Value = 0.635350030255722
$ws.Range("M14").Value = 38.45264233333334
$ws.Range("N14").Value = 115.357927
$ws.Range("O14").Value = 0.2975040117664333
$ws.Range("P14").Value = 0.2975040117664332
$ws.Range("Q14").Value = 10550.12865925705
$ws.Range("R14").Value = 94951.15793331349
$ws.Range("S14").Value = 0.1890191828770021
$ws.Range("T14").Value = 0.189019182877002
$ws.Range("G15").Value = 274.366806
$ws.Range("H15").Value = 823.100418
$ws.Range("I15").Value = 0.635350030255722
$ws.Range("J15").Value = 0.635350030255722
$ws.Range("O15").Value = 0.3694391181876273
$ws.Range("P15").Value = 0.3694391181876272
$ws.Range("Q15").Value = 13101.10141204389
$ws.Range("R15").Value = 117909.912708395
$ws.Range("S15").Value = 0.2347231549181562
$ws.Range("T15").Value = 0.2347231549181562
$ws.Range("G16").Value = 274.366806
$ws.Range("H16").Value = 823.100418
$ws.Range("I16").Value = 0.635350030255722
$ws.Range("J16").Value = 0.635350030255722
$ws.Range("M16").Value = 18.63107466666667
$ws.Range("N16").Value = 55.893224
$ws.Range("O16").Value = 0.1441466469015163
$ws.Range("P16").Value = 0.1441466469015162
$ws.Range("Q16").Value = 5111.748448640848
$ws.Range("R16").Value = 46005.73603776764
$ws.Range("S16").Value = 0.09158357647013923
$ws.Range("T16").Value = 0.09158357647013922
$ws.Range("G17").Value = 274.366806
$ws.Range("H17").Value = 823.100418
$ws.Range("I17").Value = 0.635350030255722
$ws.Range("J17").Value = 0.635350030255722
$ws.Range("M17").Value = 24.41680433333333
$ws.Range("N17").Value = 73.25041299999999
$ws.Range("O17").Value = 0.1889102231444233
$ws.Range("P17").Value = 0.1889102231444233
$ws.Range("Q17").Value = 6699.160617663625
$ws.Range("R17").Value = 60292.44555897263
$ws.Range("S17").Value = 0.1200241159904245
$ws.Range("T17").Value = 0.1200241159904245
